$wb = $excel.ActiveWorkbook

$wsStat = $wb.Worksheets.Item("StatOutput")
$wsStatMsg = $wb.Worksheets.Item("StatOutput_Message")

# Update the Cypher query text (breed filter changed from Akita to Chesapeake Bay Retriever)
$newQuery = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Chesapeake Bay Retriever']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"
$wsStatMsg.Range("A18").Value = $newQuery

# Update the resulting stat counts to reflect the new query results
$wsStat.Range("A2").Value = "0"
$wsStat.Range("B2").Value = "0"
